# Add new evidence row for team "Jimmy668" on the "Info" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Info")

$ws.Range("A3").Value = "Jimmy668"
$ws.Range("B3").Value = "iaa1tqz6k8s04pd7ps9cug85da2uaa695l5zqervd5"
$ws.Range("C3").Value = "stars1tqz6k8s04pd7ps9cug85da2uaa695l5zp85qy5"
$ws.Range("D3").Value = "juno1tqz6k8s04pd7ps9cug85da2uaa695l5zrfqxge"
$ws.Range("E3").Value = "uptick18dfa2m6jwd53tdu5e03atqtw4ysfuwytm2q08u"
$ws.Range("F3").Value = "omniflix1tqz6k8s04pd7ps9cug85da2uaa695l5zg9jycm"
$ws.Range("H3").Value = "no"
$ws.Range("G3").Value = "Jimmy668#7209"

$ws.Range("A3").Select()
